# Update "想去人数" (want-to-go count) values for the affected shows
# in both the "展览" sheet and the "全部类型" sheet.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 6003
    $ws.Range("F3").Value = 22
    $ws.Range("F5").Value = 993
    $ws.Range("F6").Value = 88
}
